# Generate Report for Handback
#
# This mirrors a localization "handback" run:
#  - the Status column for every file row flips from "Ready for handoff"
#    to "Handed back: in sync with en-US"
#  - the zh-cn and de-de sheets get their "Latest Target File" /
#    "Latest Handback File" / "Latest Handback DateTime" columns filled in
#    for both rows (a.md + b.md) now that a handback has happened
#  - a couple of long-text columns get widened so the new text isn't
#    truncated

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (C) for both data rows
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Latest Target File (I) -> hyperlink to a.md, same as column A
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/85d8ab2d5b7358e9f79f15655dd4cd528c9fbc86/e2e/a.md", "", "", "a.md")
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/85d8ab2d5b7358e9f79f15655dd4cd528c9fbc86/e2e/a.md", "", "", "a.md")
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("I3").Font.Color = 15570276

# Latest Handback File (J)
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Latest Handback DateTime (K)
$wsZh.Range("K2").Value = "2016-08-21 08:43:35"
$wsZh.Range("K3").Value = "2016-08-21 08:43:35"

# Widen Status (C) and Latest Handback File (J) columns
$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column (C) for both data rows
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Latest Target File (I) -> hyperlink to a.md, same as column A
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/85d8ab2d5b7358e9f79f15655dd4cd528c9fbc86/e2e/a.md", "", "", "a.md")
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/85d8ab2d5b7358e9f79f15655dd4cd528c9fbc86/e2e/a.md", "", "", "a.md")
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("I3").Font.Color = 15570276

# Latest Handback File (J)
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# Latest Handback DateTime (K)
$wsDe.Range("K2").Value = "2016-08-21 08:43:41"
$wsDe.Range("K3").Value = "2016-08-21 08:43:41"

# Widen Status (C) and Latest Handback File (J) columns
$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15
